$d = $word.ActiveDocument

function Merge-FollowingRunsIntoFirst($firstRunRange, $mergeEndPos) {
    # Appends all the text that currently lives, in the document, between the
    # end of $firstRunRange and $mergeEndPos onto $firstRunRange itself (via
    # InsertAfter on the *existing* run range, so its own run properties /
    # rsid attributes are preserved), then removes the now-duplicated trailing
    # runs. The net effect is that every run between the start of
    # $firstRunRange and $mergeEndPos is coalesced into a single run.
    $restStart = $firstRunRange.End
    $restText = $d.Range($restStart, $mergeEndPos).Text
    if ($restText.Length -gt 0) {
        $d.Range($restStart, $mergeEndPos).Delete()
        $firstRunRange.InsertAfter($restText)
    }
}

# --- Edit 1: in the "Infosys Ltd, Match ...Micropat..." paragraph, merge the
# leading "Infosys Ltd, " run together with the following "Match ""risk
# assessment""..." run into a single run. The "Micropat" run (and its
# surrounding proofErr spell-check markers) must stay untouched.
$rMicropat = $d.Content
if ($rMicropat.Find.Execute("Micropat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $micropatStart = $rMicropat.Start
    $anchor = $d.Range($micropatStart, $micropatStart)
    $para = $anchor.Paragraphs(1)
    $paraStart = $para.Range.Start

    # The first run of the paragraph ("Infosys Ltd, ") - find where it ends
    # by locating the comma-space right after "Infosys Ltd,".
    $rFirstRun = $d.Content
    $rFirstRun.Start = $paraStart
    $rFirstRun.End = $micropatStart
    $rFirstRun.Find.Execute("Infosys Ltd, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $firstRunEnd = $rFirstRun.End
    $firstRunRange = $d.Range($paraStart, $firstRunEnd)

    Merge-FollowingRunsIntoFirst $firstRunRange $micropatStart
}

# --- Edit 2: in the "BASF_PII_DI, Curriculum Vitae..." paragraph, merge all
# of its runs ("BASF_PII_D" + "I, " + "Curriculum Vitae...") into one run.
$rBasf = $d.Content
if ($rBasf.Find.Execute("BASF_PII_D", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $basfRunEnd = $rBasf.End
    $basfStart = $rBasf.Start
    $firstRunRange2 = $d.Range($basfStart, $basfRunEnd)

    $anchor2 = $d.Range($basfStart, $basfStart)
    $para2 = $anchor2.Paragraphs(1)
    $pEnd = $para2.Range.End
    $pText = $para2.Range.Text
    # The paragraph Range includes the trailing paragraph mark; exclude it.
    if ($pText.Length -gt 0) {
        $lastCode = [int][char]$pText.Substring($pText.Length - 1, 1)
        if ($lastCode -eq 13 -or $lastCode -eq 7) {
            $pEnd = $pEnd - 1
        }
    }

    Merge-FollowingRunsIntoFirst $firstRunRange2 $pEnd
}
